# Commit: "Remove IRA from recently edited files"
# In the BBNPPTY sheet, clear out the "banned from year X onward" flags
# (set to 1 under IRA policy) back to 0 for the "hard coal" row (row 2)
# and the "lignite" row (row 14), columns I:AE (years 2028-2050).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("BBNPPTY")

$ws.Range("I2:AE2").Value = 0
$ws.Range("I14:AE14").Value = 0

# Update the view: move/extend the selection to match the edited row
# (mirrors where the author's cursor ended up after editing row 14),
# and scroll the window so column D is left-most, matching the saved
# workbook view state.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H14:AE14").Select()
